# Apply the diff: add a new "Table_2" worksheet with capital adequacy
# ratios, and clean up a few stray empty inline-string cells on "Table_1".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Table_1")

# --- Clean up stray empty cells on Table_1 ---
$ws1.Range("B2").ClearContents()
$ws1.Range("A3").ClearContents()
$ws1.Range("B37").ClearContents()

# --- Add the new Table_2 worksheet, placed right after Table_1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Table_2"

# Values here are percentages/text labels stored as literal text in the
# source workbook, not numeric percentages - force text format first so
# Excel doesn't auto-convert "6.0%" into a numeric 0.06 cell.
$dataRange = $ws2.Range("A1:D4")
$dataRange.NumberFormat = "@"

$ws2.Range("A1").Value = "Əmsal"
$ws2.Range("B1").Value = "Norma (Sistem əhəmiyyətli)"
$ws2.Range("C1").Value = "Norma (Banklar istisna)"
$ws2.Range("D1").Value = "Fakt"

$headerRange = $ws2.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$ws2.Range("A2").Value = "9.  I dərəcəli  kapitalın  adekvatlıq əmsalı"
$ws2.Range("B2").Value = "6.0%"
$ws2.Range("C2").Value = "5.0%"
$ws2.Range("D2").Value = "11.44%"

$ws2.Range("A3").Value = "10. məcmu kapitalın  adekvatlıq  əmsalı"
$ws2.Range("B3").Value = "11.0%"
$ws2.Range("C3").Value = "9.0%"
$ws2.Range("D3").Value = "22.43%"

$ws2.Range("A4").Value = "11. Leverec əmsalı"
$ws2.Range("B4").Value = "minimum 5%"
$ws2.Range("C4").Value = "minimum 4%"
$ws2.Range("D4").Value = "4.91%"

# Keep Table_1 as the active sheet selection, matching the original workbook.
$ws1.Activate()
$ws1.Range("A1").Select() | Out-Null
